$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - Wins, Losses, Ties in AD1:AF1, matching style of other headers (s="1")
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (A1) to the new header cells so they match formatting
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2-45 with Wins=83, Losses=79, Ties=0
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 83   # AD
    $ws.Cells.Item($r, 31).Value = 79   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
